# "Generate Report for Handback"
#
# The localization-status report gets refreshed after a handback: the
# Overview/status text moves from "Ready for handoff" to
# "Handed back: in sync with en-US", and the per-language detail sheets
# (zh-cn / de-de) get their "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns populated for both rows.

$wb = $excel.ActiveWorkbook

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/5e48484403cb391a585572f481efeef969036719/e2e/62288978-c74f-438f-83ff-b02031d3c663.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/5e48484403cb391a585572f481efeef969036719/e2e/ad84fe27-7571-4181-b764-17442849d730.md"
$mdName1 = "62288978-c74f-438f-83ff-b02031d3c663.md"
$mdName2 = "ad84fe27-7571-4181-b764-17442849d730.md"
$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: status goes from "Ready for handoff" to handed-back ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.144371396019366
$wsOverview.Columns.Item(6).ColumnWidth = 29.144371396019366

# --- zh-cn detail sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column mirrors the Overview text too.
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText
$wsZh.Columns.Item(3).ColumnWidth = 29.144371396019366

# Row 2 -> 62288978 file, Row 3 -> ad84fe27 file.
$wsZh.Range("I2").Value = $mdName1
$wsZh.Range("I2").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl1, "", "", $mdName1)
$wsZh.Range("J2").Value = "62288978-c74f-438f-83ff-b02031d3c663.cf34051a15e942c44782187442fcce09126d858c.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-13 19:18:29"

$wsZh.Range("I3").Value = $mdName2
$wsZh.Range("I3").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl2, "", "", $mdName2)
$wsZh.Range("J3").Value = "ad84fe27-7571-4181-b764-17442849d730.25dd66467bfead44e6297d6567a5bd2a0349368d.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-13 19:18:29"

$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# --- de-de detail sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText
$wsDe.Columns.Item(3).ColumnWidth = 29.144371396019366

$wsDe.Range("I2").Value = $mdName1
$wsDe.Range("I2").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl1, "", "", $mdName1)
$wsDe.Range("J2").Value = "62288978-c74f-438f-83ff-b02031d3c663.cf34051a15e942c44782187442fcce09126d858c.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-13 19:18:39"

$wsDe.Range("I3").Value = $mdName2
$wsDe.Range("I3").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl2, "", "", $mdName2)
$wsDe.Range("J3").Value = "ad84fe27-7571-4181-b764-17442849d730.25dd66467bfead44e6297d6567a5bd2a0349368d.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-13 19:18:39"

$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664
